# Regenerate the "K" (strikeouts) column (column G) values for the
# miller_erik save_data sheet. The prior export used an old "Strike#"
# metric; this pass recalculates the true strikeout count (K) per
# start and writes it back into column G for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row (sheet row number) -> new K value, as recomputed upstream.
$kValues = @{
    2  = 0
    3  = 1
    4  = 2
    5  = 2
    6  = 1
    8  = 1
    9  = 1
    10 = 1
    11 = 0
    12 = 2
    13 = 2
    14 = 2
    15 = 3
    16 = 2
    17 = 1
    18 = 1
    19 = 1
    20 = 2
    21 = 2
    22 = 0
    23 = 2
    24 = 0
    25 = 1
    26 = 2
    27 = 2
    28 = 2
    29 = 1
    30 = 2
    31 = 2
    32 = 1
    33 = 1
    34 = 0
    35 = 1
    36 = 2
    37 = 1
    38 = 0
    39 = 2
    40 = 0
    41 = 0
    42 = 0
    43 = 0
    44 = 2
    45 = 2
    46 = 0
    47 = 1
    48 = 2
    49 = 0
    50 = 3
    52 = 3
    53 = 0
    54 = 1
    55 = 2
    56 = 0
    57 = 0
    58 = 2
    59 = 1
    60 = 1
    61 = 1
    62 = 1
    63 = 1
    64 = 0
    65 = 2
    66 = 1
    67 = 1
    68 = 2
    69 = 3
    70 = 1
    71 = 2
    72 = 0
    73 = 1
    74 = 1
    75 = 1
    76 = 2
    77 = 2
    78 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
